# Insert a new weekly data row at row 13 (pushing existing rows 13-106 down to
# 14-107), then populate the new row with the latest "Dulce o Americano"
# observation, continuing the existing series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 13; existing row 13 (and everything
# below it) shifts down one row.
$ws.Rows(13).Insert()

# Populate the newly inserted row 13 with the new weekly record.
$ws.Cells.Item(13, 1).Value = 7
$ws.Cells.Item(13, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(13, 3).Value = "Ñuble"
$ws.Cells.Item(13, 4).Value = 44530
$ws.Cells.Item(13, 5).Value = 16
$ws.Cells.Item(13, 6).Value = 100112024
$ws.Cells.Item(13, 7).Value = "Choclo"
$ws.Cells.Item(13, 8).Value = "Dulce o Americano"
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 10).Value = 60
$ws.Cells.Item(13, 11).Value = 16000
$ws.Cells.Item(13, 12).Value = 17000
$ws.Cells.Item(13, 13).Value = 16500
$ws.Cells.Item(13, 14).Value = "$/malla 60 unidades"
$ws.Cells.Item(13, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(13, 16).Value = 275
$ws.Cells.Item(13, 17).Value = 60
$ws.Cells.Item(13, 18).Value = "Hortaliza"
